# Weekly refresh of the "Rabanito" price series:
# - a brand-new record is inserted at row 19 (most recent date, 2022-03-04)
# - every existing record from the old row 19 down to the old row 57
#   shifts down by one row (old row N -> new row N+1)
# - the sheet grows by one row (A1:R57 -> A1:R58)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 19
$lastRow = 57
$lastCol = 18   # column R

# 1) Snapshot every existing data row (19..57) across all columns before
#    anything is overwritten.
$oldData = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $rowVals = @()
  for ($c = 1; $c -le $lastCol; $c++) {
    $rowVals += $ws.Cells.Item($r, $c).Value()
  }
  $oldData[$r] = $rowVals
}

# The brand-new last row (58) has no pre-existing cell formatting, so set
# its Fecha cell's number format (same as the rest of column D) before the
# value is written, so it doesn't fall back to a generic date format.
$ws.Cells.Item($lastRow + 1, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 2) Shift every snapshot row down by one: new row (r+1) gets old row r's
#    values, starting from the bottom so writes never clobber a value we
#    still need to read (even though we already snapshotted everything).
for ($r = $lastRow; $r -ge $firstRow; $r--) {
  $destRow = $r + 1
  $src = $oldData[$r]
  for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($destRow, $c).Value = $src[$c - 1]
  }
}

# 3) Overwrite row 19 with the new weekly record. It keeps the same
#    market/category/quality metadata as the old row 19, but with a new
#    date and a new reported volume (price fields unchanged).
$base = $oldData[$firstRow]
for ($c = 1; $c -le $lastCol; $c++) {
  $ws.Cells.Item($firstRow, $c).Value = $base[$c - 1]
}
$ws.Cells.Item($firstRow, 4).Value = 44624   # Fecha -> 2022-03-04
$ws.Cells.Item($firstRow, 10).Value = 30     # Volumen -> 30
